# Add "Area" / "Atotal" columns to the Q computation sheet (Station 6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "-" placeholders in B2/C2 become real numeric zeros so the new
# Area formulas (which multiply by B2) evaluate instead of erroring.
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# New headers.
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Per-segment area (mirrors the existing per-segment Q column E).
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Running totals, plus a small "report" pair echoing the totals.
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Match the author's final selection.
$ws.Range("J2:K2").Select() | Out-Null
